# Actualización automática 2025-08-18 16:45:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M15").Value = 605.0700000000001
$ws1.Range("O15").Value = 1.73
$ws1.Range("P15").Value = 2.12

$ws1.Range("M34").Value = "4 de 32"
$ws1.Range("O34").Value = "2 de 32"
$ws1.Range("P34").Value = "2 de 32"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F15").Value = 2509.72
$ws2.Range("F34").Value = 10337.62

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D10").Value = 183.08
$ws3.Range("E10").Value = 467.17
$ws3.Range("F10").Value = 0.2815532487504806

$ws3.Range("D16").Value = 5639.28
$ws3.Range("E16").Value = 16233.82
$ws3.Range("F16").Value = 0.257818050482099

$ws3.Range("D18").Value = 3.46
$ws3.Range("E18").Value = 1596.54
$ws3.Range("F18").Value = 0.0021625

$ws3.Range("D19").Value = 10457.7
$ws3.Range("E19").Value = 21651.58107555788
$ws3.Range("F19").Value = 0.325690879698972
